$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 314, pushing existing rows 314-334 down to 315-335.
$ws.Rows("314:314").Insert()

# Populate the newly inserted row 314 with the new weekly price record.
$ws.Cells.Item(314, 1).Value = 1
$ws.Cells.Item(314, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(314, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(314, 4).Value = 44826
$ws.Cells.Item(314, 5).Value = 15
$ws.Cells.Item(314, 6).Value = 100114013
$ws.Cells.Item(314, 7).Value = "Zanahoria"
$ws.Cells.Item(314, 8).Value = "Sin especificar"
$ws.Cells.Item(314, 9).Value = "Primera"
$ws.Cells.Item(314, 10).Value = 60
$ws.Cells.Item(314, 11).Value = 21000
$ws.Cells.Item(314, 12).Value = 22000
$ws.Cells.Item(314, 13).Value = 21500
$ws.Cells.Item(314, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(314, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(314, 16).Value = 860
$ws.Cells.Item(314, 17).Value = 25
$ws.Cells.Item(314, 18).Value = "Hortaliza"
